$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest scraped values.
# For Price cells whose new text would otherwise be auto-coerced to a
# number by Excel (losing formatting such as trailing zeros), force the
# cell to Text format first so the literal string is preserved.

$ws.Range('D2').Value = '28.694.25'
$ws.Range('E2').Value = '  +1.44%  '
$ws.Range('D3').Value = '1.560.91'
$ws.Range('E3').Value = '  -0.72%  '
$ws.Range('E4').Value = '  -0.31%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '210.40'
$ws.Range('E5').Value = '  -0.66%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.487'
$ws.Range('E6').Value = '  -0.44%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.35%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '24.78'
$ws.Range('E8').Value = '  +4.46%  '
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('E10').Value = '  -0.45%  '
$ws.Range('E11').Value = '  -0.18%  '
$ws.Range('D12').Value = '1.781.41'
$ws.Range('E12').Value = '  -0.94%  '
$ws.Range('D13').Value = '1.562.42'
$ws.Range('E13').Value = '  -0.66%  '
$ws.Range('D14').Value = '28.685.23'
$ws.Range('E14').Value = '  +1.32%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.516'
$ws.Range('E15').Value = '  +0.05%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.64'
$ws.Range('E16').Value = '  -1.09%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.51'
$ws.Range('E17').Value = '  -0.10%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '230.49'
$ws.Range('E18').Value = '  +0.33%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.38'
$ws.Range('E19').Value = '  -0.36%  '
$ws.Range('D20').Value = '0.0₃0673'
$ws.Range('E20').Value = '  -1.75%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.999'
$ws.Range('E21').Value = '  -0.22%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.92'
$ws.Range('E22').Value = '  -0.84%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.99'
$ws.Range('E23').Value = '  -0.35%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.07'
$ws.Range('E24').Value = '  +1.38%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '151.16'
$ws.Range('E25').Value = '  -0.26%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '14.80'
$ws.Range('E26').Value = '  -0.87%  '
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('E28').Value = '  -0.33%  '
$ws.Range('E29').Value = '  -1.88%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0461'
$ws.Range('E30').Value = '  -3.67%  '
$ws.Range('E31').Value = '  -1.98%  '
$ws.Range('E32').Value = '  -0.69%  '
$ws.Range('D33').Value = '1.392.30'
$ws.Range('E33').Value = '  +1.03%  '
$ws.Range('E34').Value = '  -2.59%  '
$ws.Range('E35').Value = '  -2.72%  '
$ws.Range('E36').Value = '  -1.64%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.66'
$ws.Range('E37').Value = '  +0.90%  '
$ws.Range('E38').Value = '  -3.29%  '
$ws.Range('E39').Value = '  -0.66%  '
$ws.Range('E40').Value = '  +3.84%  '
$ws.Range('E41').Value = '  -0.26%  '
$ws.Range('E42').Value = '  -0.17%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.777'
$ws.Range('E43').Value = '  -0.89%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0465'
$ws.Range('E44').Value = '  +0.14%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '64.05'
$ws.Range('E45').Value = '  +3.04%  '
$ws.Range('E46').Value = '  -2.00%  '
$ws.Range('D47').Value = '1.695.22'
$ws.Range('E47').Value = '  -0.91%  '
$ws.Range('E48').Value = '  -5.90%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '85.29'
$ws.Range('E49').Value = '  -0.21%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '43.48'
$ws.Range('E50').Value = '  +5.41%  '
$ws.Range('E51').Value = '  -0.38%  '
